$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.167.02"
$ws.Range("D3").Value = "1.858.21"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.7137"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "'240.19"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.07730"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("D9").Value = "'0.3074"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'24.92"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'0.08253"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.858.16"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'5.216"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'0.7140"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").Value = "'90.09"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "29.235.89"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "'5.862"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'243.59"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'0.000007803"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").Value = "'13.14"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "2.113.24"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'7.930"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'0.1577"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'162.54"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'8.892"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "'18.24"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.319"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.492"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'4.362"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'4.109"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "'0.05183"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("D35").Value = "'1.172"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").Value = "'0.7273"
$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").Value = "'0.01847"
$ws.Range("D39").Value = "'2.685"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "1.151.48"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'0.9019"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'6.095"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").Value = "'72.05"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'101.80"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "2.011.46"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").Value = "'0.5231"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "'1.762"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "'0.00000000119"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'9.279"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "'2.867"
$ws.Range("E51").Value = "  -0.69%  "
